$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 391.75
$ws.Range("I12").Value = 389
$ws.Range("K12").Value = 389
$ws.Range("M12").Value = -219
$ws.Range("H18").Value = 999
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H92").Value = 955
$ws.Range("J92").Value = 967.6
$ws.Range("L92").Value = 967.6
$ws.Range("N92").Value = -3463.6
$ws.Range("H94").Value = 1950
$ws.Range("I94").Value = 1950
$ws.Range("K94").Value = 1950
$ws.Range("M94").Value = -1499
$ws.Range("H96").Value = 468.4
$ws.Range("I96").Value = 514.8889
$ws.Range("J96").Value = 50
$ws.Range("K96").Value = 1544.6667
$ws.Range("L96").Value = 150
$ws.Range("M96").Value = -171.6667000000002
$ws.Range("N96").Value = -2896
$ws.Range("H100").Value = 5130.227
$ws.Range("I100").Value = 2682.7856
$ws.Range("K100").Value = 2682.7856
$ws.Range("M100").Value = -2141.7856
$ws.Range("H101").Value = 34.5
$ws.Range("J101").Value = 25
$ws.Range("L101").Value = 75
$ws.Range("N101").Value = -3319
$ws.Range("H113").Value = 8551.191999999999
$ws.Range("I113").Value = 7103.8184
$ws.Range("J113").Value = 9612.6
$ws.Range("K113").Value = 7103.8184
$ws.Range("L113").Value = 9612.6
$ws.Range("M113").Value = -3849.8184
$ws.Range("N113").Value = -16120.6
$ws.Range("H124").Value = 57956
$ws.Range("J124").Value = 57956
$ws.Range("L124").Value = 57956
$ws.Range("N124").Value = -67776
$ws.Range("H137").Value = 2344.5178
$ws.Range("I137").Value = 2001.0217
$ws.Range("K137").Value = 6003.0651
$ws.Range("M137").Value = -3453.0651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 36649.223
$ws.Range("I2").Value = 62566
$ws.Range("J2").Value = 4253.25
$ws.Range("K2").Value = 62566
$ws.Range("L2").Value = 4253.25
$ws.Range("M2").Value = -62453
$ws.Range("N2").Value = -4479.25
$ws.Range("H4").Value = 999.5
$ws.Range("I4").Value = 999.5
$ws.Range("K4").Value = 999.5
$ws.Range("M4").Value = -883.5
$ws.Range("H32").Value = 15152900
$ws.Range("I32").Value = 15874300
$ws.Range("K32").Value = 15874300
$ws.Range("M32").Value = -15874013
$ws.Range("H69").Value = 79999.164
$ws.Range("J69").Value = 79999.164
$ws.Range("L69").Value = 79999.164
$ws.Range("N69").Value = -81497.164
$ws.Range("H72").Value = 79999.164
$ws.Range("J72").Value = 79999.164
$ws.Range("L72").Value = 239997.492
$ws.Range("N72").Value = -247485.492
$ws.Range("H74").Value = 1509.4286
$ws.Range("I74").Value = 1531.0264
$ws.Range("J74").Value = 1304.25
$ws.Range("K74").Value = 1531.0264
$ws.Range("L74").Value = 1304.25
$ws.Range("M74").Value = -657.0264
$ws.Range("N74").Value = -3052.25
$ws.Range("H77").Value = 1509.4286
$ws.Range("I77").Value = 1531.0264
$ws.Range("J77").Value = 1304.25
$ws.Range("K77").Value = 7655.132
$ws.Range("L77").Value = 6521.25
$ws.Range("M77").Value = -3287.132
$ws.Range("N77").Value = -15257.25
$ws.Range("H95").Value = 99471
$ws.Range("J95").Value = 99471
$ws.Range("L95").Value = 99471
$ws.Range("N95").Value = -104963
$ws.Range("H116").Value = 36649.223
$ws.Range("I116").Value = 62566
$ws.Range("J116").Value = 4253.25
$ws.Range("K116").Value = 62566
$ws.Range("L116").Value = 4253.25
$ws.Range("M116").Value = -60272
$ws.Range("N116").Value = -8841.25
$ws.Range("H132").Value = 2809.5417
$ws.Range("I132").Value = 2164.9546
$ws.Range("J132").Value = 9900
$ws.Range("K132").Value = 6494.8638
$ws.Range("L132").Value = 29700
$ws.Range("M132").Value = -3964.8638
$ws.Range("N132").Value = -34760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 36649.223
$ws.Range("I3").Value = 62566
$ws.Range("J3").Value = 4253.25
$ws.Range("K3").Value = 62566
$ws.Range("L3").Value = 4253.25
$ws.Range("M3").Value = -62452
$ws.Range("N3").Value = -4481.25
$ws.Range("H58").Value = 34621.5
$ws.Range("I58").Value = 51743
$ws.Range("K58").Value = 51743
$ws.Range("M58").Value = -51449
$ws.Range("H70").Value = 79999.164
$ws.Range("J70").Value = 79999.164
$ws.Range("L70").Value = 79999.164
$ws.Range("N70").Value = -80585.164
$ws.Range("H73").Value = 79999.164
$ws.Range("J73").Value = 79999.164
$ws.Range("L73").Value = 79999.164
$ws.Range("N73").Value = -82027.164
$ws.Range("H94").Value = 2109
$ws.Range("I94").Value = 1984.4286
$ws.Range("J94").Value = 2377.3076
$ws.Range("K94").Value = 1984.4286
$ws.Range("L94").Value = 2377.3076
$ws.Range("M94").Value = -1533.4286
$ws.Range("N94").Value = -3279.3076
$ws.Range("H134").Value = 2017
$ws.Range("I134").Value = 1584.4348
$ws.Range("J134").Value = 5333.3335
$ws.Range("K134").Value = 4753.3044
$ws.Range("L134").Value = 16000.0005
$ws.Range("M134").Value = -2218.3044
$ws.Range("N134").Value = -21070.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 335
$ws.Range("I25").Value = 335
$ws.Range("K25").Value = 335
$ws.Range("M25").Value = -161
$ws.Range("H31").Value = 1130.909
$ws.Range("I31").Value = 1171.7142
$ws.Range("K31").Value = 1171.7142
$ws.Range("M31").Value = -876.7141999999999
$ws.Range("H34").Value = 1130.909
$ws.Range("I34").Value = 1171.7142
$ws.Range("K34").Value = 1171.7142
$ws.Range("M34").Value = -969.7141999999999
$ws.Range("H58").Value = 1576.258
$ws.Range("I58").Value = 1265.1364
$ws.Range("J58").Value = 2336.7778
$ws.Range("K58").Value = 1265.1364
$ws.Range("L58").Value = 2336.7778
$ws.Range("M58").Value = -1062.1364
$ws.Range("N58").Value = -2742.7778
$ws.Range("H86").Value = 25019.45
$ws.Range("I86").Value = 35475
$ws.Range("J86").Value = 14563.9
$ws.Range("K86").Value = 35475
$ws.Range("L86").Value = 14563.9
$ws.Range("M86").Value = -34352
$ws.Range("N86").Value = -16809.9
$ws.Range("H89").Value = 25019.45
$ws.Range("I89").Value = 35475
$ws.Range("J89").Value = 14563.9
$ws.Range("K89").Value = 177375
$ws.Range("L89").Value = 72819.5
$ws.Range("M89").Value = -171759
$ws.Range("N89").Value = -84051.5
$ws.Range("H94").Value = 10638.272
$ws.Range("J94").Value = 1983.7142
$ws.Range("L94").Value = 1983.7142
$ws.Range("N94").Value = -2885.7142
$ws.Range("H136").Value = 1576.258
$ws.Range("I136").Value = 1265.1364
$ws.Range("J136").Value = 2336.7778
$ws.Range("K136").Value = 3795.4092
$ws.Range("L136").Value = 7010.3334
$ws.Range("M136").Value = -1245.4092
$ws.Range("N136").Value = -12110.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7372.25
$ws.Range("I56").Value = 7372.25
$ws.Range("K56").Value = 7372.25
$ws.Range("M56").Value = -6842.25
$ws.Range("H94").Value = 2356.1428
$ws.Range("I94").Value = 298.8
$ws.Range("J94").Value = 7499.5
$ws.Range("K94").Value = 896.4000000000001
$ws.Range("L94").Value = 22498.5
$ws.Range("M94").Value = -220.4000000000001
$ws.Range("N94").Value = -23850.5
$ws.Range("H121").Value = 1587.9412
$ws.Range("J121").Value = 2137.7273
$ws.Range("L121").Value = 6413.1819
$ws.Range("N121").Value = -9033.1819
$ws.Range("H122").Value = 8014.2
$ws.Range("J122").Value = 9761.625
$ws.Range("L122").Value = 87854.625
$ws.Range("N122").Value = -92754.625
$ws.Range("H129").Value = 1183
$ws.Range("I129").Value = 675.7143
$ws.Range("J129").Value = 2366.6667
$ws.Range("K129").Value = 2027.1429
$ws.Range("L129").Value = 7100.000100000001
$ws.Range("M129").Value = 2972.8571
$ws.Range("N129").Value = -17100.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5344.737
$ws.Range("I132").Value = 5475
$ws.Range("K132").Value = 16425
$ws.Range("M132").Value = -13895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7796.875
$ws.Range("I7").Value = 11300
$ws.Range("K7").Value = 11300
$ws.Range("M7").Value = -11188
$ws.Range("H46").Value = 1882.6
$ws.Range("I46").Value = 557
$ws.Range("K46").Value = 557
$ws.Range("M46").Value = -369
$ws.Range("H93").Value = 3085.0278
$ws.Range("I93").Value = 1954.3704
$ws.Range("J93").Value = 6477
$ws.Range("K93").Value = 1954.3704
$ws.Range("L93").Value = 6477
$ws.Range("M93").Value = -706.3704
$ws.Range("N93").Value = -8973
$ws.Range("H126").Value = 7796.875
$ws.Range("I126").Value = 11300
$ws.Range("K126").Value = 33900
$ws.Range("M126").Value = -31430

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4381.6665
$ws.Range("I96").Value = 4190.6
$ws.Range("J96").Value = 4477.2
$ws.Range("K96").Value = 4190.6
$ws.Range("L96").Value = 4477.2
$ws.Range("M96").Value = -2817.6
$ws.Range("N96").Value = -7223.2
$ws.Range("H126").Value = 2921.625
$ws.Range("I126").Value = 2728.1667
$ws.Range("J126").Value = 3502
$ws.Range("K126").Value = 8184.500100000001
$ws.Range("L126").Value = 10506
$ws.Range("M126").Value = -5714.500100000001
$ws.Range("N126").Value = -15446
